$wb = $excel.ActiveWorkbook

# Rename Sheet2 -> Order, Sheet3 -> Address
$wsOrder = $wb.Worksheets.Item(2)
$wsOrder.Name = "Order"
$wsAddress = $wb.Worksheets.Item(3)
$wsAddress.Name = "Address"

# Customer sheet: add Orders (list) and Address (class) columns
$wsCustomer = $wb.Worksheets.Item(1)
$wsCustomer.Range("F1").Value = "Orders"
$wsCustomer.Range("F2").Value = "1,2"
$wsCustomer.Range("G1").Value = "Address"
$wsCustomer.Range("G2").Value = 1
[void]$wsCustomer.Range("G3").Select()

# Order sheet: populate with Id/Name rows
$wsOrder.Range("A1").Value = "Id"
$wsOrder.Range("B1").Value = "Name"
$wsOrder.Range("A2").Value = 1
$wsOrder.Range("B2").Value = "A"
$wsOrder.Range("A3").Value = 2
$wsOrder.Range("B3").Value = "B"
$wsOrder.Range("A4").Value = 3
$wsOrder.Range("B4").Value = "C"
$wsOrder.Range("A5").Value = 4
$wsOrder.Range("B5").Value = "D"
[void]$wsOrder.Range("B5").Select()

# Address sheet: populate with Id/Name rows
$wsAddress.Range("A1").Value = "Id"
$wsAddress.Range("B1").Value = "Name"
$wsAddress.Range("A2").Value = 1
$wsAddress.Range("B2").Value = "A"
$wsAddress.Range("A3").Value = 2
$wsAddress.Range("B3").Value = "B"
[void]$wsAddress.Range("B3").Select()

# Return focus to the Customer sheet so it remains the selected/active tab
[void]$wsCustomer.Select()
[void]$wsCustomer.Range("G3").Select()

Write-Host "Edit complete"
